$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.314.98"
$ws.Range("E2").Value = "  +3.65%  "
$ws.Range("D3").Value = "3.122.30"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'219.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.52%  "
$ws.Range("D6").Value = "'625.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("D7").Value = "'0.386"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.01%  "
$ws.Range("D8").Value = "'0.959"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +20.16%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "3.117.72"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").Value = "'0.723"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +20.56%  "
$ws.Range("E12").Value = "  +5.90%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.25%  "
$ws.Range("D14").Value = "'34.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.43%  "
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").Value = "'5.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "91.180.40"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "3.698.72"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "3.108.89"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").Value = "'3.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +15.78%  "
$ws.Range("E20").Value = "  +10.27%  "
$ws.Range("D21").Value = "'14.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.61%  "
$ws.Range("D22").Value = "'436.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("D23").Value = "'8.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.52%  "
$ws.Range("D24").Value = "'5.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.13%  "
$ws.Range("E25").Value = "  +12.05%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'12.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'86.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.34%  "
$ws.Range("D28").Value = "3.290.17"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").Value = "'9.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.84%  "
$ws.Range("E32").Value = "  -7.44%  "
$ws.Range("D33").Value = "'528.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.27%  "
$ws.Range("E34").Value = "  +5.08%  "
$ws.Range("D35").Value = "'7.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.86%  "
$ws.Range("E36").Value = "  +8.85%  "
$ws.Range("D37").Value = "'23.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.40%  "
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'0.0835"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +23.24%  "
$ws.Range("D43").Value = "'0.150"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.34%  "
$ws.Range("D45").Value = "'0.379"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("E46").Value = "  +5.63%  "
$ws.Range("D47").Value = "'147.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'44.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").Value = "'1.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.93%  "
$ws.Range("D50").Value = "'166.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.90%  "
$ws.Range("E51").Value = "  +21.13%  "
